# Apply the two logical changes described by the diff:
#  1. "sets" sheet: D4 (home_points for set 3 of match 1) changes 4 -> 5
#  2. "rallies" sheet: append a new row 68 with rally data (dimension grows
#     from A1:P67 to A1:P68)

$wb = $excel.ActiveWorkbook

# --- 1. sets!D4 4 -> 5 -----------------------------------------------------
$setsWs = $wb.Worksheets.Item("sets")
$setsWs.Range("D4").Value = 5

# --- 2. rallies!A68:P68 new row --------------------------------------------
$ralliesWs = $wb.Worksheets.Item("rallies")

$ralliesWs.Range("A68").Value = 67
$ralliesWs.Range("B68").Value = 1
$ralliesWs.Range("C68").Value = 3
$ralliesWs.Range("D68").Value = 5
$ralliesWs.Range("E68").Value = "NOS"
$ralliesWs.Range("F68").Value = ""
$ralliesWs.Range("G68").Value = 3
$ralliesWs.Range("H68").Value = "LINHA"
$ralliesWs.Range("I68").Value = "PONTO"
$ralliesWs.Range("J68").Value = "NOS"
$ralliesWs.Range("K68").Value = 5
$ralliesWs.Range("L68").Value = 0
$ralliesWs.Range("M68").Value = "1 3 l"
$ralliesWs.Range("N68").Value = "FRENTE"
$ralliesWs.Range("O68").Value = "FRENTE"
$ralliesWs.Range("P68").Value = "FRENTE"
